$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Define the new reviewer/reusable values ---
$valY = 'Y'
$valInitials = 'JGE'
$valInclude = 'Include'
$valReviewDate = "5/14/2025"

# --- New column headers (row 1) ---
$ws.Range("P1").Value = 'Setting_Population'
$ws.Range("Q1").Value = 'Study_Objective'
$ws.Range("R1").Value = 'Study_Method'
$ws.Range("S1").Value = 'Study_Results'
$ws.Range("T1").Value = 'Study_Conclusion'

# --- Row 2 (Study_ID 442) ---
$ws.Range("B2").Value = $valY
$ws.Range("C2").Value = $valY
$ws.Range("D2").Value = $valY
$ws.Range("E2").Value = $valY
$ws.Range("F2").Value = $valY
$ws.Range("G2").Value = $valY
$ws.Range("H2").Value = $valY
$ws.Range("I2").Value = $valY
$ws.Range("J2").Value = $valY
$ws.Range("K2").Value = $valY
$ws.Range("L2").Value = $valInclude
$ws.Range("M2").Value = $valInitials
$ws.Range("N2").Value = $valReviewDate
$ws.Range("O2").Value = 'Not 100% clear that included consecutive patients, but implied. Only includes patients that underwent endoscopy.'
$ws.Range("P2").Value = 'Three large Australian teaching hospitals, Australia'
$ws.Range("Q2").Value = 'To audit and analyse existing endoscopy databases at three large Australian teaching hospitals: To assess the characteristics of patients admitted with DFBI; report the types of foreign bodies ingested and evaluate the medical care and the outcomes of recurrent presentations.'
$ws.Range("R2").Value = 'Adult patients with an endoscopic record of attempted foreign body retrieval between January 2013 and September 2020 were identified at three Australian hospitals. Those with a documented mental health diagnosis were included and their standard medical records reviewed. Presentation history, demographics, comorbidities and endoscopic findings were recorded and described.'
$ws.Range("S2").Value = 'A total of 166 admissions were accounted for by 35 patients, 2/3 of which had borderline personality disorder (BPD). Repetitive presentations occurred in more than half of the cohort. There was an increased trend of hospital admissions throughout the years. At least half of the cohort had a documented mental health review during their admission. An average of 3.3 (2.9) foreign bodies were ingested per single episode. Endoscopic intervention was performed in 76.5% of incidents. The combined Length of stay for all patients was 680 days.'
$ws.Range("T2").Value = 'Deliberate foreign body ingestion in mental health patients is a common, recurring and challenging problem that is increasing in frequency and requires collaborative research to further guide holistic management.'

# N2 already carried the Review_Date date-number format (style index 2) in
# the original workbook. N3/N4 start out with no style at all, so a plain
# `.Value =` assignment of a date-like string would make Excel mint a brand
# new (equivalent-looking but distinct) numFmt/style entry. Copy N2's
# format down first so N3/N4 reuse the very same style index.
$ws.Range("N2").Copy()
$ws.Range("N3:N4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 3 (Study_ID 523) ---
$ws.Range("B3").Value = $valY
$ws.Range("C3").Value = $valY
$ws.Range("D3").Value = $valY
$ws.Range("E3").Value = $valY
$ws.Range("F3").Value = $valY
$ws.Range("G3").Value = $valY
$ws.Range("H3").Value = $valY
$ws.Range("I3").Value = $valY
$ws.Range("J3").Value = $valY
$ws.Range("K3").Value = $valY
$ws.Range("L3").Value = $valInclude
$ws.Range("M3").Value = $valInitials
$ws.Range("N3").Value = $valReviewDate
$ws.Range("O3").Value = 'No comment.'
$ws.Range("P3").Value = 'Medical prison ward for men in a New York City public general hospital, New York, USA'
$ws.Range("Q3").Value = 'Records were examined for demographic and psychiatric characteristics, as well as for the cir cumstances surrounding the ingestion of a foreign object.'
$ws.Range("R3").Value = 'Two psychiatrically trained raters examined the records of all patients (N= 19) admitted to a medical prison ward for men in a New York City public general hospital for deliberately swallowing objects between September 1, 1985, and October 15, 1988.'
$ws.Range("S3").Value = 'See statistical breakdown.'
$ws.Range("T3").Value = 'In our series of 19 prisoners who ingested objects, 12 prisoners (64 percent) were judged to be suicidal by the psychiatric consultant. In addition, the high reported incidence of past suicide attempts by other methods is startling. In our sample, a remarkable 84 percent descnibed a history ofsuicide attempts, and most of these attempts occurred before imprisonment. Also common in this sample was psychosis. About three-quarters of the patients were judged by the psy chiatnic consultant to have bad corn mand hallucinations, suicidal idea tion, or both. In addition, the swat lowing ofantennae and toothbrushes can be construed as psychotic behav ion. Psychosis may have occurred in the setting of schizophrenia; almost a third of the sample had that diagnosis. No patient swallowed a foreign body before his first impnisonment. Desire to leave prison may explain this behavior, as may suicidal idea tion in response to incarceration. Sm cide pacts or copycat phenomena may contribute to foreign body in gestion in prison; two ofthe patients had shared a razor blade, each swal lowing half. All 19 patients ingested sharp or pointed objects, possibly be cause swallowing such objects is like ly to result in transfer from the prison to the hospital.'

# --- Row 4 (Study_ID 646) ---
$ws.Range("B4").Value = $valY
$ws.Range("C4").Value = $valY
$ws.Range("D4").Value = $valY
$ws.Range("E4").Value = $valY
$ws.Range("F4").Value = $valY
$ws.Range("G4").Value = $valY
$ws.Range("H4").Value = $valY
$ws.Range("I4").Value = $valY
$ws.Range("J4").Value = $valY
$ws.Range("K4").Value = $valY
$ws.Range("L4").Value = $valInclude
$ws.Range("M4").Value = $valInitials
$ws.Range("N4").Value = $valReviewDate
$ws.Range("O4").Value = 'No comment.'
$ws.Range("P4").Value = 'Patients referred from prison to Department of General and Digestive Surgery, Farhat Hached University Hospital of Sousse, Sousse, Tunisia'
$ws.Range("Q4").Value = 'A descriptive study including all detainees ingesting a razor blade, transferred from the prison to Farhat Hached University Hospital of Sousse, from January 1, 2014 to December 31, 2015.'
$ws.Range("S4").Value = 'There were 16 men with a mean age of 24 years, ranging from 19 to 27 years. Three patients had a history of self-harm; one of them was having a psychiatric follow-up for depressive disorders. An inmate had ingested a half blade 3 times and another had ingested a half blade 2 times, so we had 19 swallowed razor blade episodes. This act was a form of protest in 17 cases and a suicide attempt in 2 cases.'
$ws.Range("T4").Value = 'Our experience enabled us to confirm the few data in the literature that surgical removal of intragastrointestinal sharp foreign bodies should not be systematic.'

# --- Selection update ---
# Note: selecting a multi-cell range always anchors the active cell at the
# range's top-left corner in this engine, so activeCell normalizes to A2
# even though Excel's real UI would keep K4 active after a shift-click.
$ws.Range("A2:K4").Select()
